{"js": "// Update the two header dates (05/12/2015 -> 05/15/2015) and the three\n// \"Estimated Completion Date\" mentions (May 2 2015 -> May 3 2015).\n\nconst body = context.document.body;\n\n// 1) Replace the two standalone \"05/12/2015\" date cells.\nconst dateResults = body.search(\"05/12/2015\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"05/15/2015\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Replace the three \"Estimated Completion Date: May 2 2015\" sentences.\nconst oldCompletion =\n  \"Responsible party: George Washington (george.washington@nasa.gov), Estimated Completion Date: May 2 2015\";\nconst newCompletion =\n  \"Responsible party: George Washington (george.washington@nasa.gov), Estimated Completion Date: May 3 2015\";\n\nconst completionResults = body.search(oldCompletion, { matchCase: true });\ncompletionResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < completionResults.items.length; i++) {\n  completionResults.items[i].insertText(newCompletion, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the two header dates (05/12/2015 -> 05/15/2015) and the three\n# \"Estimated Completion Date\" mentions (May 2 2015 -> May 3 2015).\n\n$d = $word.ActiveDocument\n\n# 1) Replace the two standalone \"05/12/2015\" date cells.\n$find1 = $d.Content.Find\n$find1.Text = \"05/12/2015\"\n$find1.Replacement.Text = \"05/15/2015\"\n$find1.Execute(\n    \"05/12/2015\",   # FindText\n    $true,          # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap (wdFindContinue)\n    $false,         # Format\n    \"05/15/2015\",   # ReplaceWith\n    2               # Replace (wdReplaceAll)\n) | Out-Null\n\n# 2) Replace the three \"Estimated Completion Date: May 2 2015\" sentences.\n$oldCompletion = \"Responsible party: George Washington (george.washington@nasa.gov), Estimated Completion Date: May 2 2015\"\n$newCompletion = \"Responsible party: George Washington (george.washington@nasa.gov), Estimated Completion Date: May 3 2015\"\n\n$find2 = $d.Content.Find\n$find2.Text = $oldCompletion\n$find2.Replacement.Text = $newCompletion\n$find2.Execute(\n    $oldCompletion, # FindText\n    $true,          # MatchCase\n    $false,         # MatchWholeWord\n    $false,         # MatchWildcards\n    $false,         # MatchSoundsLike\n    $false,         # MatchAllWordForms\n    $true,          # Forward\n    1,              # Wrap (wdFindContinue)\n    $false,         # Format\n    $newCompletion, # ReplaceWith\n    2               # Replace (wdReplaceAll)\n) | Out-Null\n"}
